# "not db yet, but made a bunch of other improvements"
# Rework the "school" sheet from a 2-column key/value dump
# (link / school+Fremont / email) into a proper 4-column table:
#   School Name | background_link | buffer_amount | users_name
#   Fremont     | <wolf image url>| 70            | Karrie

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("school")

# Start clean - the old layout (A1 "link" header, A2/B2 "school"/"Fremont",
# A3 "email") is being fully replaced.
$ws.Cells.Clear()

# Populate in the same order the values were authored in, so the shared
# string table lines up with how Excel would have built it incrementally.
$ws.Range("A2").Value = "Fremont"
$ws.Range("B1").Value = "background_link"
$ws.Range("A1").Value = "School Name"
$ws.Range("D1").Value = "users_name"
$ws.Range("D2").Value = "Karrie"

# Bold header row
$ws.Range("A1:D1").Font.Bold = $true

# Hyperlink cell with the background image link
$ws.Hyperlinks.Add($ws.Range("B2"), "https://static.vecteezy.com/system/resources/previews/038/035/644/large_2x/ai-generated-wolf-high-quality-image-free-photo.jpg")

# Buffer amount column + value
$ws.Range("C1").Value = "buffer_amount"
$ws.Range("C2").Value = 70

# Roomier columns now that B/C hold a long URL and numbers
$ws.Range("A1").EntireColumn.ColumnWidth = 11.6
$ws.Range("B1").EntireColumn.ColumnWidth = 42.4
$ws.Range("C1").EntireColumn.ColumnWidth = 42.4
$ws.Range("D1").EntireColumn.ColumnWidth = 10.8

# Leave the cursor where the author's session ended up
$ws.Range("C3").Select()
